$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 4400
$ws.Range("J12").Value = 4400
$ws.Range("L12").Value = 4400
$ws.Range("N12").Value = -4740

$ws.Range("H70").Value = 100004250
$ws.Range("I70").Value = 1985
$ws.Range("J70").Value = 142862350
$ws.Range("K70").Value = 5955
$ws.Range("L70").Value = 428587050
$ws.Range("M70").Value = -5685
$ws.Range("N70").Value = -428587590

$ws.Range("H73").Value = 100004250
$ws.Range("I73").Value = 1985
$ws.Range("J73").Value = 142862350
$ws.Range("K73").Value = 5955
$ws.Range("L73").Value = 428587050
$ws.Range("M73").Value = -5019
$ws.Range("N73").Value = -428588922

$ws.Range("H80").Value = 1073.3158
$ws.Range("J80").Value = 1095.1875
$ws.Range("L80").Value = 3285.5625
$ws.Range("N80").Value = -5281.5625

$ws.Range("H83").Value = 1073.3158
$ws.Range("J83").Value = 1095.1875
$ws.Range("L83").Value = 9856.6875
$ws.Range("N83").Value = -19840.6875

$ws.Range("H111").Value = 3799.889
$ws.Range("I111").Value = 3799.889
$ws.Range("K111").Value = 11399.667
$ws.Range("M111").Value = -8332.667000000001

$ws.Range("H113").Value = 4002.762
$ws.Range("I113").Value = 3847.4375
$ws.Range("K113").Value = 3847.4375
$ws.Range("M113").Value = -593.4375

$ws.Range("H137").Value = 1823.6666
$ws.Range("I137").Value = 1695.2
$ws.Range("K137").Value = 5085.6
$ws.Range("M137").Value = -2535.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 86.666664
$ws.Range("I4").Value = 90
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 90
$ws.Range("L4").Value = 80
$ws.Range("M4").Value = 26
$ws.Range("N4").Value = -312

$ws.Range("H39").Value = 1002500
$ws.Range("I39").Value = 1002500
$ws.Range("K39").Value = 1002500
$ws.Range("M39").Value = -1001980

$ws.Range("H132").Value = 1426.9062
$ws.Range("I132").Value = 1426.9062
$ws.Range("K132").Value = 4280.7186
$ws.Range("M132").Value = -1750.7186

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1118.5862
$ws.Range("I20").Value = 809.64703
$ws.Range("J20").Value = 1556.25
$ws.Range("K20").Value = 809.64703
$ws.Range("L20").Value = 1556.25
$ws.Range("M20").Value = -562.64703
$ws.Range("N20").Value = -2050.25

$ws.Range("H35").Value = 57057.2
$ws.Range("J35").Value = 57057.2
$ws.Range("L35").Value = 57057.2
$ws.Range("N35").Value = -57677.2

$ws.Range("H38").Value = 7036
$ws.Range("J38").Value = 7036
$ws.Range("L38").Value = 7036
$ws.Range("N38").Value = -7868

$ws.Range("H49").Value = 24000
$ws.Range("J49").Value = 24000
$ws.Range("L49").Value = 24000
$ws.Range("N49").Value = -24478

$ws.Range("H99").Value = 2941.8
$ws.Range("I99").Value = 1291.6666
$ws.Range("K99").Value = 1291.6666
$ws.Range("M99").Value = 206.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 5166.5
$ws.Range("I39").Value = 5166.5
$ws.Range("K39").Value = 5166.5
$ws.Range("M39").Value = -4775.5

$ws.Range("H49").Value = 5166.5
$ws.Range("I49").Value = 5166.5
$ws.Range("K49").Value = 5166.5
$ws.Range("M49").Value = -4984.5

$ws.Range("H58").Value = 1743.6072
$ws.Range("I58").Value = 1551.0476
$ws.Range("J58").Value = 2321.2856
$ws.Range("K58").Value = 1551.0476
$ws.Range("L58").Value = 2321.2856
$ws.Range("M58").Value = -1348.0476
$ws.Range("N58").Value = -2727.2856

$ws.Range("H134").Value = 2366.5334
$ws.Range("I134").Value = 2366.5334
$ws.Range("K134").Value = 7099.600199999999
$ws.Range("M134").Value = -4564.600199999999

$ws.Range("H136").Value = 1743.6072
$ws.Range("I136").Value = 1551.0476
$ws.Range("J136").Value = 2321.2856
$ws.Range("K136").Value = 4653.142800000001
$ws.Range("L136").Value = 6963.8568
$ws.Range("M136").Value = -2103.142800000001
$ws.Range("N136").Value = -12063.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 705.44446
$ws.Range("J97").Value = 659.8
$ws.Range("L97").Value = 1979.4
$ws.Range("N97").Value = -2971.4

$ws.Range("H137").Value = 3086.25
$ws.Range("I137").Value = 1442.1666
$ws.Range("J137").Value = 4072.7
$ws.Range("K137").Value = 4326.4998
$ws.Range("L137").Value = 12218.1
$ws.Range("M137").Value = 773.5002000000004
$ws.Range("N137").Value = -22418.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 23772.6
$ws.Range("J15").Value = 23772.6
$ws.Range("L15").Value = 23772.6
$ws.Range("N15").Value = -24348.6

$ws.Range("H17").Value = 135.8
$ws.Range("J17").Value = 135.8
$ws.Range("L17").Value = 135.8
$ws.Range("N17").Value = -471.8

$ws.Range("H81").Value = 23772.6
$ws.Range("J81").Value = 23772.6
$ws.Range("L81").Value = 23772.6
$ws.Range("N81").Value = -25768.6

$ws.Range("H84").Value = 23772.6
$ws.Range("J84").Value = 23772.6
$ws.Range("L84").Value = 71317.79999999999
$ws.Range("N84").Value = -81301.79999999999

$ws.Range("H132").Value = 2275.7144
$ws.Range("I132").Value = 1928.5
$ws.Range("K132").Value = 5785.5
$ws.Range("M132").Value = -3255.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 950
$ws.Range("I19").Value = 950
$ws.Range("K19").Value = 950
$ws.Range("M19").Value = -780

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 500
$ws.Range("I25").Value = 500
$ws.Range("K25").Value = 500
$ws.Range("M25").Value = -207

$ws.Range("H62").Value = 20427.285
$ws.Range("J62").Value = 21230
$ws.Range("L62").Value = 21230
$ws.Range("N62").Value = -22478

$ws.Range("H65").Value = 20427.285
$ws.Range("J65").Value = 21230
$ws.Range("L65").Value = 106150
$ws.Range("N65").Value = -112390

$ws.Range("H81").Value = 4922.5
$ws.Range("I81").Value = 4922.5
$ws.Range("K81").Value = 9845
$ws.Range("M81").Value = -8784

$ws.Range("H84").Value = 4922.5
$ws.Range("I84").Value = 4922.5
$ws.Range("K84").Value = 49225
$ws.Range("M84").Value = -43921

$ws.Range("H103").Value = 30288.445
$ws.Range("J103").Value = 30288.445
$ws.Range("L103").Value = 30288.445
$ws.Range("N103").Value = -32632.445

$ws.Range("H132").Value = 3789.6667
$ws.Range("I132").Value = 3476.8823
$ws.Range("K132").Value = 10430.6469
$ws.Range("M132").Value = -7900.6469
